# NEXT AuthorList update: remove Hafidi and Meziani from the author list.
#
# In the original sheet, row 45 is "Hafidi" and row 66 is "Meziani".
# Deleting row 45 first shifts every row below it up by one, so the row
# that used to be 66 (Meziani) becomes row 65 before it is removed.
# This mirrors selecting each author's row and choosing Delete in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove "Hafidi" row.
$ws.Rows(45).Select()
$ws.Rows(45).Delete()

# Remove "Meziani" row (shifted up to 65 after the previous delete).
$ws.Rows(65).Select()
$ws.Rows(65).Delete()
